$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("strategy_id-0")

# --- Part 1: reset calibration scalars (rows 4-10,12-15,21-22) back to 1 across J:AS ---
$ws.Range("J4:AS4").Value = 1
$ws.Range("J5:AS5").Value = 1
$ws.Range("J6:AS6").Value = 1
$ws.Range("J7:AS7").Value = 1
$ws.Range("J8:AS8").Value = 1
$ws.Range("J9:AS9").Value = 1
$ws.Range("J10:AS10").Value = 1
$ws.Range("J12:AS12").Value = 1
$ws.Range("J13:AS13").Value = 1
$ws.Range("J14:AS14").Value = 1
$ws.Range("J15:AS15").Value = 1
$ws.Range("J21:AS21").Value = 1
$ws.Range("J22:AS22").Value = 1

# --- Part 2: update IPPU production values (rows 96-104, 111-112) with new per-column figures ---
$row96 = New-Object 'object[,]' 1,36
$row96[0,0] = 3540248.08932218
$row96[0,1] = 1381388.08405699
$row96[0,2] = 947872.303277281
$row96[0,3] = 1128940.26152003
$row96[0,4] = 1182372.61931889
$row96[0,5] = 840378.76471821
$row96[0,6] = 840378.8
$row96[0,7] = 855311.104238914
$row96[0,8] = 870508.733721496
$row96[0,9] = 885976.4028899251
$row96[0,10] = 901718.909955133
$row96[0,11] = 917741.1383852561
$row96[0,12] = 934048.058420526
$row96[0,13] = 950644.728615089
$row96[0,14] = 967536.29740621
$row96[0,15] = 984728.004711371
$row96[0,16] = 1002225.18355373
$row96[0,17] = 1020033.26171649
$row96[0,18] = 1038157.76342663
$row96[0,19] = 1056604.31106857
$row96[0,20] = 1075378.6269283
$row96[0,21] = 1094486.53496849
$row96[0,22] = 1113933.96263509
$row96[0,23] = 1133726.94269613
$row96[0,24] = 1153871.61511313
$row96[0,25] = 1174374.22894575
$row96[0,26] = 1195241.1442903
$row96[0,27] = 1216478.83425275
$row96[0,28] = 1238093.88695667
$row96[0,29] = 1260093.00758701
$row96[0,30] = 1282483.02047003
$row96[0,31] = 1305270.87119032
$row96[0,32] = 1328463.62874538
$row96[0,33] = 1352068.48773844
$row96[0,34] = 1376092.77061035
$row96[0,35] = 1400543.92991104
$ws.Range("J96:AS96").Value = $row96

$row97 = New-Object 'object[,]' 1,36
$row97[0,0] = 236567.586256145
$row97[0,1] = 204116.588128651
$row97[0,2] = 235223.687703097
$row97[0,3] = 298351.441278218
$row97[0,4] = 254931.648160624
$row97[0,5] = 206840.297546668
$row97[0,6] = 206840.3
$row97[0,7] = 216892.280242579
$row97[0,8] = 227432.764450765
$row97[0,9] = 238485.492834809
$row97[0,10] = 250075.359326577
$row97[0,11] = 262228.46764786
$row97[0,12] = 274972.190103485
$row97[0,13] = 288335.229231638
$row97[0,14] = 302347.682450261
$row97[0,15] = 317041.109845115
$row97[0,16] = 332448.605252192
$row97[0,17] = 348604.870794584
$row97[0,18] = 365546.295041668
$row97[0,19] = 383311.034966658
$row97[0,20] = 401939.101887115
$row97[0,21] = 421472.451581973
$row97[0,22] = 441955.078788051
$row97[0,23] = 463433.116288891
$row97[0,24] = 485954.938819089
$row97[0,25] = 509571.272018148
$row97[0,26] = 534335.306679248
$row97[0,27] = 560302.818550254
$row97[0,28] = 587532.29395678
$row97[0,29] = 616085.06153027
$row97[0,30] = 646025.430337753
$row97[0,31] = 677420.834724417
$row97[0,32] = 710341.986195196
$row97[0,33] = 744863.032677475
$row97[0,34] = 781061.7255236059
$row97[0,35] = 819019.594629376
$ws.Range("J97:AS97").Value = $row97

$row98 = New-Object 'object[,]' 1,36
$row98[0,0] = 22990.5828986593
$row98[0,1] = 20426.7556584196
$row98[0,2] = 20588.0726635464
$row98[0,3] = 19782.3715166274
$row98[0,4] = 17569.1677205607
$row98[0,5] = 18633.2100072617
$row98[0,6] = 18633.21
$row98[0,7] = 18955.4050699888
$row98[0,8] = 19283.171357343
$row98[0,9] = 19616.6051964446
$row98[0,10] = 19955.8045874358
$row98[0,11] = 20300.8692250228
$row98[0,12] = 20651.9005277769
$row98[0,13] = 21009.0016679427
$row98[0,14] = 21372.2776017618
$row98[0,15] = 21741.8351003206
$row98[0,16] = 22117.7827809314
$row98[0,17] = 22500.2311390564
$row98[0,18] = 22889.292580784
$row98[0,19] = 23285.0814558657
$row98[0,20] = 23687.7140913251
$row98[0,21] = 24097.3088256478
$row98[0,22] = 24513.9860435626
$row98[0,23] = 24937.8682114237
$row98[0,24] = 25369.0799132051
$row98[0,25] = 25807.7478871177
$row98[0,26] = 26254.0010628584
$row98[0,27] = 26707.9705995048
$row98[0,28] = 27169.7899240638
$row98[0,29] = 27639.5947706879
$row98[0,30] = 28117.5232205688
$row98[0,31] = 28603.7157425209
$row98[0,32] = 29098.3152342672
$row98[0,33] = 29601.4670644383
$row98[0,34] = 30113.3191152979
$row98[0,35] = 30634.0218262076
$ws.Range("J98:AS98").Value = $row98

$row99 = New-Object 'object[,]' 1,36
$row99[0,0] = 5943.71776390242
$row99[0,1] = 5582.97987596999
$row99[0,2] = 6115.48269826848
$row99[0,3] = 6409.13054586766
$row99[0,4] = 6600.71488356777
$row99[0,5] = 6893.55115918622
$row99[0,6] = 6893.551
$row99[0,7] = 6986.18228846381
$row99[0,8] = 7080.0582990748
$row99[0,9] = 7175.19575764183
$row99[0,10] = 7271.61161472483
$row99[0,11] = 7369.32304865495
$row99[0,12] = 7468.34746859511
$row99[0,13] = 7568.70251764188
$row99[0,14] = 7670.40607596884
$row99[0,15] = 7773.47626401237
$row99[0,16] = 7877.93144570005
$row99[0,17] = 7983.79023172263
$row99[0,18] = 8091.07148284984
$row99[0,19] = 8199.7943132908
$row99[0,20] = 8309.97809409961
$row99[0,21] = 8421.642456626631
$row99[0,22] = 8534.807296016201
$row99[0,23] = 8649.492774751359
$row99[0,24] = 8765.719326246181
$row99[0,25] = 8883.507658486311
$row99[0,26] = 9002.878757718579
$row99[0,27] = 9123.85389219004
$row99[0,28] = 9246.454615937369
$row99[0,29] = 9370.702772627061
$row99[0,30] = 9496.620499447359
$row99[0,31] = 9624.230231052399
$row99[0,32] = 9753.554703559341
$row99[0,33] = 9884.616958599299
$row99[0,34] = 10017.4403474226
$row99[0,35] = 10152.0485350593
$ws.Range("J99:AS99").Value = $row99

$row100 = New-Object 'object[,]' 1,36
$row100[0,0] = 763830.428372727
$row100[0,1] = 745504.701947103
$row100[0,2] = 898053.4390704751
$row100[0,3] = 418933.677949677
$row100[0,4] = 704274.564686539
$row100[0,5] = 1045524.07351848
$row100[0,6] = 1045524
$row100[0,7] = 1284958.83124427
$row100[0,8] = 1579226.49120693
$row100[0,9] = 1940884.21347692
$row100[0,10] = 2385364.95626093
$row100[0,11] = 2931635.97037283
$row100[0,12] = 3603008.18548778
$row100[0,13] = 4428130.95346248
$row100[0,14] = 5442214.59723327
$row100[0,15] = 6688532.9349124
$row100[0,16] = 8220269.89603666
$row100[0,17] = 10102789.030308
$row100[0,18] = 12416422.7551849
$row100[0,19] = 15259900.3674111
$row100[0,20] = 18754561.101431
$row100[0,21] = 23049532.017815
$row100[0,22] = 28328091.6768423
$row100[0,23] = 34815491.1531976
$row100[0,24] = 42788566.1366051
$row100[0,25] = 52587550.2939293
$row100[0,26] = 64630594.0023247
$row100[0,27] = 79431608.0088549
$row100[0,28] = 97622193.45601919
$row100[0,29] = 119978594.089421
$row100[0,30] = 147454820.774534
$row100[0,31] = 181223361.839401
$row100[0,32] = 222725216.468788
$row100[0,33] = 273731386.216252
$row100[0,34] = 336418448.650924
$row100[0,35] = 413461438.080332
$ws.Range("J100:AS100").Value = $row100

$row101 = New-Object 'object[,]' 1,36
$row101[0,0] = 128643.470396428
$row101[0,1] = 129961.584459341
$row101[0,2] = 146626.26763033
$row101[0,3] = 159722.989456785
$row101[0,4] = 156231.572344347
$row101[0,5] = 166113.674968843
$row101[0,6] = 166113.7
$row101[0,7] = 174694.184535095
$row101[0,8] = 183717.887870669
$row101[0,9] = 193207.704157313
$row101[0,10] = 203187.710126617
$row101[0,11] = 213683.226176545
$row101[0,12] = 224720.880612134
$row101[0,13] = 236328.677204501
$row101[0,14] = 248536.066239558
$row101[0,15] = 261374.019236704
$row101[0,16] = 274875.107527051
$row101[0,17] = 289073.584890559
$row101[0,18] = 304005.474461727
$row101[0,19] = 319708.660124337
$row101[0,20] = 336222.982627135
$row101[0,21] = 353590.340664285
$row101[0,22] = 371854.797177078
$row101[0,23] = 391062.691146564
$row101[0,24] = 411262.755160765
$row101[0,25] = 432506.239054735
$row101[0,26] = 454847.039937151
$row101[0,27] = 478341.838933349
$row101[0,28] = 503050.24499171
$row101[0,29] = 529034.946118271
$row101[0,30] = 556361.868423252
$row101[0,31] = 585100.343383009
$row101[0,32] = 615323.283741793
$row101[0,33] = 647107.368499586
$row101[0,34] = 680533.23745535
$row101[0,35] = 715685.695799267
$ws.Range("J101:AS101").Value = $row101

$row102 = New-Object 'object[,]' 1,36
$row102[0,0] = 5656387.83548591
$row102[0,1] = 3593915.00889202
$row102[0,2] = 2672859.80298353
$row102[0,3] = 2317597.84976661
$row102[0,4] = 2136732.98874682
$row102[0,5] = 2648716.73710727
$row102[0,6] = 2648717
$row102[0,7] = 2843299.8916517
$row102[0,8] = 3052177.44057465
$row102[0,9] = 3276399.77622663
$row102[0,10] = 3517094.17380295
$row102[0,11] = 3775470.72159945
$row102[0,12] = 4052828.40471741
$row102[0,13] = 4350561.63569606
$row102[0,14] = 4670167.26490548
$row102[0,15] = 5013252.10594452
$row102[0,16] = 5381541.01387755
$row102[0,17] = 5776885.5569233
$row102[0,18] = 6201273.32519265
$row102[0,19] = 6656837.92327486
$row102[0,20] = 7145869.69690999
$row102[0,21] = 7670827.24767551
$row102[0,22] = 8234349.7935773
$row102[0,23] = 8839270.43768758
$row102[0,24] = 9488630.4115378
$row102[0,25] = 10185694.3648749
$row102[0,26] = 10933966.7786503
$row102[0,27] = 11737209.5837568
$row102[0,28] = 12599461.0740932
$row102[0,29] = 13525056.2090397
$row102[0,30] = 14518648.4074161
$row102[0,31] = 15585232.9424907
$row102[0,32] = 16730172.0556595
$row102[0,33] = 17959221.9150521
$row102[0,34] = 19278561.5546004
$row102[0,35] = 20694823.9390603
$ws.Range("J102:AS102").Value = $row102

$row103 = New-Object 'object[,]' 1,36
$row103[0,0] = 93123.94161505481
$row103[0,1] = 91798.3360602827
$row103[0,2] = 86492.44993581819
$row103[0,3] = 82729.3955594726
$row103[0,4] = 85686.298610992
$row103[0,5] = 96845.46375741941
$row103[0,6] = 96845.46000000001
$row103[0,7] = 100347.399581237
$row103[0,8] = 103975.969577886
$row103[0,9] = 107735.7489559
$row103[0,10] = 111631.482257005
$row103[0,11] = 115668.085585937
$row103[0,12] = 119850.652814172
$row103[0,13] = 124184.462007984
$row103[0,14] = 128674.98208895
$row103[0,15] = 133327.879735285
$row103[0,16] = 138149.026532743
$row103[0,17] = 143144.506384089
$row103[0,18] = 148320.623186499
$row103[0,19] = 153683.908786572
$row103[0,20] = 159241.131223007
$row103[0,21] = 164999.303267322
$row103[0,22] = 170965.691273412
$row103[0,23] = 177147.824347114
$row103[0,24] = 183553.503847332
$row103[0,25] = 190190.813230732
$row103[0,26] = 197068.12825242
$row103[0,27] = 204194.127535478
$row103[0,28] = 211577.803522691
$row103[0,29] = 219228.473824296
$row103[0,30] = 227155.79297606
$row103[0,31] = 235369.764622538
$row103[0,32] = 243880.75414087
$row103[0,33] = 252699.501721064
$row103[0,34] = 261837.135919258
$row103[0,35] = 271305.187701069
$ws.Range("J103:AS103").Value = $row103

$row104 = New-Object 'object[,]' 1,36
$row104[0,0] = 39043.7759426143
$row104[0,1] = 38173.1957319621
$row104[0,2] = 34409.9105640365
$row104[0,3] = 37487.2158903309
$row104[0,4] = 37748.2858024121
$row104[0,5] = 40936.2795187436
$row104[0,6] = 40936.28
$row104[0,7] = 43133.9629078672
$row104[0,8] = 45449.6294274239
$row104[0,9] = 47889.6135628057
$row104[0,10] = 50460.5893620563
$row104[0,11] = 53169.5891725397
$row104[0,12] = 56024.0228764037
$row104[0,13] = 59031.6981587067
$row104[0,14] = 62200.8418636493
$row104[0,15] = 65540.12249732421
$row104[0,16] = 69058.6739385371
$row104[0,17] = 72766.1204225531
$row104[0,18] = 76672.60286610819
$row104[0,19] = 80788.8066056892
$row104[0,20] = 85125.99062495799
$row104[0,21] = 89696.0183512618
$row104[0,22] = 94511.3901054689
$row104[0,23] = 99585.2772938886
$row104[0,24] = 104931.5584358
$row104[0,25] = 110564.857125134
$row104[0,26] = 116500.58203015
$row104[0,27] = 122754.969040504
$row104[0,28] = 129345.125677016
$row104[0,29] = 136289.077885578
$row104[0,30] = 143605.819343231
$row104[0,31] = 151315.363411253
$row104[0,32] = 159438.797877371
$row104[0,33] = 167998.342636836
$row104[0,34] = 177017.410470137
$row104[0,35] = 186520.671083586
$ws.Range("J104:AS104").Value = $row104

$row111 = New-Object 'object[,]' 1,36
$row111[0,0] = 1425.35643314755
$row111[0,1] = 1423.0691481715
$row111[0,2] = 1233.91172926136
$row111[0,3] = 1104.89494825993
$row111[0,4] = 1164.77810124642
$row111[0,5] = 1402.04092335277
$row111[0,6] = 1402.041
$row111[0,7] = 1487.1143207011
$row111[0,8] = 1577.34973715768
$row111[0,9] = 1673.0604760355
$row111[0,10] = 1774.57877002983
$row111[0,11] = 1882.25701111701
$row111[0,12] = 1996.46897378332
$row111[0,13] = 2117.61111247717
$row111[0,14] = 2246.10393778825
$row111[0,15] = 2382.3934761309
$row111[0,16] = 2526.95281799828
$row111[0,17] = 2680.28376016196
$row111[0,18] = 2842.91854751711
$row111[0,19] = 3015.42172061977
$row111[0,20] = 3198.39207532932
$row111[0,21] = 3392.46474135859
$row111[0,22] = 3598.31338694653
$row111[0,23] = 3816.65255730657
$row111[0,24] = 4048.2401549677
$row111[0,25] = 4293.88007061826
$row111[0,26] = 4554.42497358454
$row111[0,27] = 4830.77927163064
$row111[0,28] = 5123.90225035355
$row111[0,29] = 5434.81140307119
$row111[0,30] = 5764.585962762
$row111[0,31] = 6114.3706483162
$row111[0,32] = 6485.37963810293
$row111[0,33] = 6878.90078464622
$row111[0,34] = 7296.30008503989
$row111[0,35] = 7739.02642261921
$ws.Range("J111:AS111").Value = $row111

$row112 = New-Object 'object[,]' 1,36
$row112[0,0] = 8171.54536731342
$row112[0,1] = 6849.26936407516
$row112[0,2] = 6168.47833038071
$row112[0,3] = 6150.69402598503
$row112[0,4] = 6379.92732824989
$row112[0,5] = 6281.74226743705
$row112[0,6] = 6281.742
$row112[0,7] = 6797.25384193988
$row112[0,8] = 7355.07121937934
$row112[0,9] = 7958.66594070047
$row112[0,10] = 8611.794728618641
$row112[0,11] = 9318.522601708361
$row112[0,12] = 10083.2481747365
$row112[0,13] = 10910.7310352703
$row112[0,14] = 11806.1213669495
$row112[0,15] = 12774.9920037954
$row112[0,16] = 13823.373115061
$row112[0,17] = 14957.7897364961
$row112[0,18] = 16185.3023816208
$row112[0,19] = 17513.5509857664
$row112[0,20] = 18950.8024563904
$row112[0,21] = 20506.0021256116
$row112[0,22] = 22188.8294252041
$row112[0,23] = 24009.7581305659
$row112[0,24] = 25980.1215486144
$row112[0,25] = 28112.1830553347
$row112[0,26] = 30419.2124219987
$row112[0,27] = 32915.5684051042
$row112[0,28] = 35616.788114066
$row112[0,29] = 38539.6837128768
$row112[0,30] = 41702.4470575995
$row112[0,31] = 45124.7629209481
$row112[0,32] = 48827.9315086549
$row112[0,33] = 52835.0010301572
$row112[0,34] = 57170.9111487122
$row112[0,35] = 61862.6482037606
$ws.Range("J112:AS112").Value = $row112

